# Update "Arbeidsmarkedstilknytning for Vestfold kommuner og hele landet 2023 SSB 13563.xlsx"
#
# The lower "Prosent" table (row 17 header) previously re-used the same long
# column headers as the "I alt" table in row 5. The edit shortens the header
# labels for the percentage table (D17:H17) to more compact wording, fixes
# cell C23 to use an explicit formula instead of the shared formula, and
# moves the active selection to I15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bosatte")

# Shortened header labels for the "Prosent" (percentage) table's header row.
$ws.Range("D17").Value = "Arbeidsledige"
$ws.Range("E17").Value = "Arbeidsmarkedstiltak"
$ws.Range("F17").Value = "Utdanning"
$ws.Range("G17").Value = "AAP / uføretrygd"
$ws.Range("H17").Value = " AFP/alderspensjon"

# C23 ("Tønsberg" row, Sysselsatte %) — re-enter as its own formula.
$ws.Range("C23").Formula = "=(C11/B11)"

# Move the active selection to I15, matching the saved cursor position.
$ws.Range("I15").Select()
